$d = $word.ActiveDocument

# --- 1. After the "composer install ... lalu enter" paragraph, add a note run ---
$pComposer = $d.Paragraphs.Item(16)
$rComposer = $d.Range($pComposer.Range.Start, $pComposer.Range.End - 1)
$rComposer.InsertAfter("(pastikan computer Anda telah terinstall composer)")

# --- 2. After the "php spark serve ... lalu tekan enter" paragraph, add a note run ---
$d2 = $word.ActiveDocument
$pServe = $d2.Paragraphs.Item(17)
$rServe = $d2.Range($pServe.Range.Start, $pServe.Range.End - 1)
$rServe.InsertAfter("(terminal jangan di close)")

# --- 3. Add <w:lastRenderedPageBreak/> inside the run that starts the "Fitur edit" heading ---
$d3 = $word.ActiveDocument
$pHeading = $d3.Paragraphs.Item(26)
$rHeading = $d3.Range($pHeading.Range.Start, $pHeading.Range.End - 1)
$headingXml = '<w:p w14:paraId="2A7820C9" w14:textId="396F53E0" w:rsidR="0052371A" w:rsidRPr="0052371A" w:rsidRDefault="0052371A" w:rsidP="0052371A"><w:pPr><w:pStyle w:val="Heading2"/></w:pPr><w:r><w:lastRenderedPageBreak/><w:t>Fitur edit</w:t></w:r></w:p>'
$rHeading.InsertXML($headingXml)

# --- 4. Remove <w:lastRenderedPageBreak/> from the run holding the screenshot picture ---
$d4 = $word.ActiveDocument
$pPicture = $d4.Paragraphs.Item(29)
$rPicture = $d4.Range($pPicture.Range.Start, $pPicture.Range.End)
$pictureXml = '<w:p w14:paraId="3AD70F71" w14:textId="1AFDB42D" w:rsidR="0052371A" w:rsidRDefault="0052371A" w:rsidP="0052371A"><w:pPr><w:jc w:val="center"/></w:pPr><w:r w:rsidRPr="0052371A"><w:rPr><w:noProof/></w:rPr><w:drawing><wp:inline distT="0" distB="0" distL="0" distR="0" wp14:anchorId="7DA29248" wp14:editId="435124CD"><wp:extent cx="3590925" cy="3003947"/><wp:effectExtent l="0" t="0" r="0" b="6350"/><wp:docPr id="3" name="Picture 3"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="1" name=""/><pic:cNvPicPr/></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId7"/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr><a:xfrm><a:off x="0" y="0"/><a:ext cx="3592545" cy="3005302"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom></pic:spPr></pic:pic></a:graphicData></a:graphic></wp:inline></w:drawing></w:r></w:p>'
$rPicture.InsertXML($pictureXml)

Write-Output "edits applied"
